$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.833.49'
$ws.Range("E2").Value = '  +2.28%  '
$ws.Range("D3").Value = '2.122.80'
$ws.Range("E3").Value = '  +10.52%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '334.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.58%  '
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5239'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4413'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +8.45%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09065'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +8.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '46.86'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +10.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.189'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +7.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '25.43'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.99%  '
$ws.Range("D13").Value = '2.121.13'
$ws.Range("E13").Value = '  +10.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.782'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.894'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +8.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '98.45'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001141'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.10%  '
$ws.Range("E18").Value = '  -0.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06639'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.96%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.412'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.70%  '
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").Value = '30.930.15'
$ws.Range("E23").Value = '  +2.56%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.78%  '
$ws.Range("D25").Value = '2.368.23'
$ws.Range("E25").Value = '  +10.87%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.255'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.03'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.586'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +14.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '163.51'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.52%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.94'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.83%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.182'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.40%  '
$ws.Range("E32").Value = '  +2.45%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.265'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.036'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.40%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.573'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +31.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02603'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.597'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.02%  '
$ws.Range("B38").Value = 'FraxShare'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '9.596'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +11.75%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06767'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.75'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +11.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2275'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.74%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6849'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.262'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.25%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.13'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6435'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.25%  '
$ws.Range("E46").Value = '  -0.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.267'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.64%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.677'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.288'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '83.30'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07083'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.58%  '
